## Alteração de fundo da pagina de abertura do Pitch
## Troca o destaque (highlight) amarelo pelo vermelho no titulo da
## disciplina, na primeira pagina (slide de abertura) do Pitch.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("CaixaDeTexto 16")

$tr = $sh.TextFrame.TextRange

# Seleciona todo o texto da caixa (intervalo explícito de caracteres) e
# troca o realce amarelo (FFFF00) por vermelho (FF0000) em todos os
# trechos de texto ("DISCIPLINA:", "PROJETO DE SISTEMAS ..." e
# "QUALIDADE DE SOFTWARE ...").
$full = $tr.Characters(1, $tr.Length)
$full.Font.Highlight.RGB = 255   # 255 = RGB(255,0,0) -> vermelho (FF0000)
